$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 0.75
$ws.Range("C5").Value = 0.75
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 626
$ws.Range("F5").Value = 627
$ws.Range("G5").Value = 0.22
$ws.Range("H5").Value = 0.01
$ws.Range("I5").Value = 0.47
$ws.Range("J5").Value = 0.01

$ws.Range("L5").Select()
